$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 225, pushing the existing rows 225:235 down to 226:236
$ws.Rows("225:225").Insert()

# Populate the newly inserted row 225 with the new weekly record
$ws.Range("A225").Value = 4
$ws.Range("B225").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C225").Value = "Los Lagos"
$ws.Range("D225").Value = 44753
$ws.Range("E225").Value = 10
$ws.Range("F225").Value = 100112039
$ws.Range("G225").Value = "Ciboulette"
$ws.Range("H225").Value = "Sin especificar"
$ws.Range("I225").Value = "Primera"
$ws.Range("J225").Value = 80
$ws.Range("K225").Value = 3500
$ws.Range("L225").Value = 3500
$ws.Range("M225").Value = 3500
$ws.Range("N225").Value = "`$/docena de atados"
$ws.Range("O225").Value = "Región Metropolitana"
$ws.Range("P225").Value = 1167
$ws.Range("Q225").Value = 3
$ws.Range("R225").Value = "Hortaliza"
